$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helpers
# ------------------------------------------------------------------

# Finds findText starting the search at character position $searchFrom
# (to the end of the document) and overwrites it in place with
# $replaceText, using a plain Range.Text assignment (NOT Find's built-in
# Replace) so that Word does not renormalise/merge unrelated runs that
# happen to live in the same paragraph.
# Returns the end position of the replacement (for chaining repeated
# searches, e.g. when the same text occurs more than once), or -1 if
# not found.
function ReplaceNext($searchFrom, $findText, $replaceText) {
    $r = $d.Range($searchFrom, $d.Content.End)
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = $replaceText
        return $r.End
    }
    return -1
}

# Inserts $text at $pos (collapsed range), optionally italicising it,
# and returns the position right after the inserted text.
function InsertPiece($pos, $text, $italic) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $len = $text.Length
    if ($italic) {
        $rr = $d.Range($pos, $pos + $len)
        $rr.Font.Italic = 1
    }
    return $pos + $len
}

# ------------------------------------------------------------------
# 1. Title: "Circle Language Spec: Execution Control" -> "... Execution Flow"
# ------------------------------------------------------------------
ReplaceNext 0 "Circle Language Spec: Execution Control" "Circle Language Spec: Execution Flow" | Out-Null

# ------------------------------------------------------------------
# 2. "Concept" paragraph: rewrite with wiggle-room wording + italics + bookmark
#    "Jumps are a form of execution control. Execution control is explained
#     in the article Execution Control."
#    ->
#    "Jumps are a form of control over execution flow. Execution flow is
#     explained in the article Execution Flow." (+ _GoBack bookmark before
#     the final ".")
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Jumps are a form of execution control. Execution control is explained in the article Execution Control.")
if ($found) {
    $start = $rng.Start
    $rng.Text = ""

    $pos = $start
    $pos = InsertPiece $pos "Jumps are a form of control over " $false
    $pos = InsertPiece $pos "execution flow" $true
    $pos = InsertPiece $pos ". Execution flow is explained in the article " $false
    $pos = InsertPiece $pos "Execution Flow" $true
    $pos = InsertPiece $pos "." $false

    if ($d.Bookmarks.Exists("_GoBack")) {
        $bm = $d.Bookmarks.Item("_GoBack")
        $bm.Delete()
    }
    $bmRange = $d.Range($pos - 1, $pos - 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# 3. Remaining plain "execution control" -> "execution flow" occurrences,
#    one at a time (in document order), each a surgical Range.Text
#    overwrite so unrelated sibling runs in the same paragraph are left
#    untouched.
# ------------------------------------------------------------------
$pos = 0

$pos = ReplaceNext $pos "Jumps are a form of execution control explained in the article" "Jumps are a form of execution flow explained in the article"

$pos = ReplaceNext $pos "order is also indicated for calls to execution control statements" "order is also indicated for calls to execution flow statements"

$pos = ReplaceNext $pos "Normal execution order is a form of execution control. Following the " "Normal execution order is a form of execution flow. Following the "

$pos = ReplaceNext $pos "are execution control statements explained in the article " "are execution flow statements explained in the article "

$pos = ReplaceNext $pos "order is also indicated for calls to execution control statements" "order is also indicated for calls to execution flow statements"

$pos = ReplaceNext $pos "statement is an execution control statement, that immediately jumps to another part of the code." "statement is an execution flow statement, that immediately jumps to another part of the code."

$pos = ReplaceNext $pos "is an execution control command that is passed a reference to the clause to go to. The reference to the clause to go to is called the " "is an execution flow command that is passed a reference to the clause to go to. The reference to the clause to go to is called the "

$pos = ReplaceNext $pos "The implementation of the execution control command " "The implementation of the execution flow command "

$pos = ReplaceNext $pos " must somehow first roll back part of the call stack, so that the call to the execution control command does not return to where it was called from. Next the " " must somehow first roll back part of the call stack, so that the call to the execution flow command does not return to where it was called from. Next the "

$pos = ReplaceNext $pos ". A jump is a type of execution control statement, explained by the article " ". A jump is a type of execution flow statement, explained by the article "

$pos = ReplaceNext $pos "The execution control command " "The execution flow command "

$pos = ReplaceNext $pos "The implementation of the execution control command " "The implementation of the execution flow command "
